$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E185:E234").Value = "M2"

$ws.Range("E186:E234").Select()
